$d = $word.ActiveDocument

$replacements = @(
    @("2025-02-11 Tuesday", "2025-02-12 Wednesday"),
    @("828×2=1656", "930×3=2790"),
    @("192×9=1728", "350×4=1400"),
    @("906×2=1812", "990×2=1980"),
    @("542×2=1084", "926×6=5556"),
    @("155×3=465", "605×4=2420"),
    @("230×8=1840", "280×9=2520"),
    @("164×2=328", "587×4=2348"),
    @("130×5=650", "110×5=550"),
    @("443×2=886", "674×5=3370"),
    @("289×6=1734", "837×3=2511"),
    @("542×7=3794", "499×2=998"),
    @("930×6=5580", "584×7=4088"),
    @("212×8=1696", "631×7=4417"),
    @("167×7=1169", "889×6=5334"),
    @("112×8=896", "476×7=3332"),
    @("527×2=1054", "558×4=2232"),
    @("197×6=1182", "363×8=2904"),
    @("940×3=2820", "774×6=4644"),
    @("647×3=1941", "135×6=810"),
    @("240×6=1440", "600×4=2400"),
    @("970×9=8730", "977×3=2931"),
    @("121×2=242", "525×6=3150"),
    @("980×6=5880", "517×8=4136"),
    @("637×4=2548", "702×2=1404"),
    @("596×5=2980", "188×7=1316")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
